$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column B first so the shared-string table is rebuilt fresh in row order,
# matching the target layout exactly (old stale entries get garbage collected).
$ws.Range("B2:B67").ClearContents()

$rows = @(
    @{B='East Fall - North Fall'; D=4.73514010152673; E=9.70908955318368; G=0.487701763959324; H=0.999998144928142},
    @{B='East Fall - West Fall'; D=-5.5444454025375; E=9.08727180497897; G=-0.610133109422309; H=0.999981079972281},
    @{B='East Spring - East Fall'; D=-12.7193707110061; E=0.958838523872745; G=-13.2653939055688; H=0},
    @{B='East Spring - East Summer'; D=-13.8408691026889; E=0.967187008513982; G=-14.3104373620098; H=0},
    @{B='East Spring - North Fall'; D=-7.98423060947941; E=9.7043607151421; G=-0.822746685108407; H=0.999630192664193},
    @{B='East Spring - North Spring'; D=-6.96785313737891; E=9.7440730074095; G=-0.715086302419992; H=0.999906197776597},
    @{B='East Spring - North Summer'; D=-19.5258949921872; E=9.68463771404193; G=-2.0161719590064; H=0.682430848991085},
    @{B='East Spring - West Fall'; D=-18.2638161135436; E=9.08219730846464; G=-2.01094685495565; H=0.686054313903753},
    @{B='East Spring - West Spring'; D=0.66500807133076; E=9.17868842728819; G=0.0724513177017421; H=1},
    @{B='East Spring - West Summer'; D=-12.6019688598669; E=9.04113089682697; G=-1.39384873459686; H=0.964898744210969},
    @{B='East Summer - East Fall'; D=1.12149839168281; E=0.994741398395681; G=1.12742708154257; H=0.993494483418537},
    @{B='East Summer - North Fall'; D=5.85663849320954; E=9.69841038102244; G=0.603876126408265; H=0.99998297316374},
    @{B='East Summer - North Summer'; D=-5.68502588949826; E=9.67865342467764; G=-0.587377772511532; H=0.999987184054399},
    @{B='East Summer - West Fall'; D=-4.42294701085468; E=9.0756949099743; G=-0.487339763481231; H=0.999998159332012},
    @{B='East Summer - West Summer'; D=1.23890024282206; E=9.03450703912146; G=0.137129811007655; H=0.999999999998002},
    @{B='East Winter - East Fall'; D=-17.8745995915453; E=0.970139778635657; G=-18.4247672193; H=0},
    @{B='East Winter - East Spring'; D=-5.15522888053912; E=0.939241157249013; G=-5.48871697194202; H=0.00000284877470302369},
    @{B='East Winter - East Summer'; D=-18.9960979832281; E=0.984585974459425; G=-19.2934883047239; H=0},
    @{B='East Winter - North Fall'; D=-13.1394594900185; E=9.70726366805329; G=-1.35356985648393; H=0.971770282390022},
    @{B='East Winter - North Spring'; D=-12.123082017918; E=9.74699795076949; G=-1.24377598919685; H=0.985366293858215},
    @{B='East Winter - North Summer'; D=-24.6811238727263; E=9.68755117844326; G=-2.54771545647643; H=0.310577426485677},
    @{B='East Winter - North Winter'; D=-1.96775395662565; E=9.78924461576805; G=-0.201011828170694; H=0.999999999866597},
    @{B='East Winter - West Fall'; D=-23.4190449940827; E=9.08532944253276; G=-2.57767702780782; H=0.292779025434662},
    @{B='East Winter - West Spring'; D=-4.49022080920836; E=9.18179587556852; G=-0.48903513757654; H=0.999998090998738},
    @{B='East Winter - West Summer'; D=-17.757197740406; E=9.04429662106161; G=-1.96335862084117; H=0.718429756520325},
    @{B='East Winter - West Winter'; D=6.92482794038639; E=9.19983626300987; G=0.752712085564969; H=0.999844281809665},
    @{B='North Spring - East Fall'; D=-5.75151757362722; E=9.74880694719963; G=-0.589971429814738; H=0.99998659065913},
    @{B='North Spring - East Summer'; D=-6.87301596531003; E=9.73798649867825; G=-0.705794361723845; H=0.999917646155434},
    @{B='North Spring - North Fall'; D=-1.0163774721005; E=1.79177978526806; G=-0.56724463600779; H=0.99999105183554},
    @{B='North Spring - North Summer'; D=-12.5580418548083; E=1.74437922175264; G=-7.19914666387209; H=0},
    @{B='North Spring - West Fall'; D=-11.2959629761647; E=10.4510411946597; G=-1.08084570386506; H=0.995471661973648},
    @{B='North Spring - West Summer'; D=-5.63411572248798; E=10.4153357128104; G=-0.540944226652075; H=0.999994526667027},
    @{B='North Summer - East Fall'; D=6.80652428118107; E=9.68937949289619; G=0.702472669810415; H=0.999921429718469},
    @{B='North Summer - North Fall'; D=11.5416643827078; E=1.31131556667712; G=8.80159183342456; H=0},
    @{B='North Summer - West Fall'; D=1.26207887864358; E=10.3956622042315; G=0.121404375579832; H=0.999999999999624},
    @{B='North Winter - East Fall'; D=-15.9068456349196; E=9.79105145914457; G=-1.62463099099158; H=0.900253280954528},
    @{B='North Winter - East Spring'; D=-3.18747492391346; E=9.78635254000651; G=-0.325706120935568; H=0.999999974853457},
    @{B='North Winter - East Summer'; D=-17.0283440266024; E=9.78038834120294; G=-1.74107033714246; H=0.849043412527023},
    @{B='North Winter - North Fall'; D=-11.1717055333929; E=1.87236539558874; G=-5.96662679181811; H=0.000000174912725858789},
    @{B='North Winter - North Spring'; D=-10.1553280612924; E=2.01915095472738; G=-5.0295041277206; H=0.0000333754844398237},
    @{B='North Winter - North Summer'; D=-22.7133699161007; E=1.84767907464409; G=-12.2929193861634; H=0},
    @{B='North Winter - West Fall'; D=-21.4512910374571; E=10.4904811661359; G=-2.04483385439969; H=0.662344605749412},
    @{B='North Winter - West Spring'; D=-2.5224668525827; E=10.5741592787432; G=-0.238550109383497; H=0.999999999137878},
    @{B='North Winter - West Summer'; D=-15.7894437837804; E=10.4550183922808; G=-1.51022630389996; H=0.93810141108339},
    @{B='West Fall - North Fall'; D=10.2795855040642; E=10.4140408084054; G=0.987089036156587; H=0.997973910423427},
    @{B='West Spring - East Fall'; D=-13.3843787823369; E=9.18371553649476; G=-1.45740345823532; H=0.951625297123289},
    @{B='West Spring - East Summer'; D=-14.5058771740197; E=9.17221538534404; G=-1.58150202155066; H=0.915948283470481},
    @{B='West Spring - North Fall'; D=-8.64923868081017; E=10.4983590272866; G=-0.823865773529911; H=0.99962533120852},
    @{B='West Spring - North Spring'; D=-7.63286120870967; E=10.5349863391915; G=-0.724525021956076; H=0.999893159668856},
    @{B='West Spring - North Summer'; D=-20.190903063518; E=10.4801178301737; G=-1.92659122642548; H=0.742564710058987},
    @{B='West Spring - West Fall'; D=-18.9288241848744; E=1.85958293866404; G=-10.1790696135732; H=0},
    @{B='West Spring - West Summer'; D=-13.2669769311977; E=1.85639507726336; G=-7.1466344064839; H=0},
    @{B='West Summer - East Fall'; D=-0.117401851139242; E=9.04624238363353; G=-0.0129779687698447; H=1},
    @{B='West Summer - North Fall'; D=4.61773825038748; E=10.3783890964127; G=0.444937861501416; H=0.999999294847626},
    @{B='West Summer - North Summer'; D=-6.92392613232031; E=10.3599227167666; G=-0.668337623900861; H=0.999952308726854},
    @{B='West Summer - West Fall'; D=-5.66184725367674; E=1.30834677535969; G=-4.3274821020751; H=0.000944535817102743},
    @{B='West Winter - East Fall'; D=-24.7994275319316; E=9.20175205971669; G=-2.69507669528483; H=0.228966901340645},
    @{B='West Winter - East Spring'; D=-12.0800568209255; E=9.1967345498796; G=-1.31351587407551; H=0.977538995733244},
    @{B='West Winter - East Summer'; D=-25.9209259236145; E=9.19027250880143; G=-2.82047413706071; H=0.171714581059956},
    @{B='West Winter - North Fall'; D=-20.0642874304049; E=10.5141430315297; G=-1.90831410322614; H=0.754239528094991},
    @{B='West Winter - North Spring'; D=-19.0479099583044; E=10.5507121781958; G=-1.80536722418313; H=0.8153680302627},
    @{B='West Winter - North Summer'; D=-31.6059518131127; E=10.495928807785; G=-3.01125821181924; H=0.105641010592493},
    @{B='West Winter - North Winter'; D=-8.89258189701205; E=10.5898289642568; G=-0.839728566629985; H=0.999550219414321},
    @{B='West Winter - West Fall'; D=-30.3438729344691; E=1.94497882011735; G=-15.6011328352862; H=0},
    @{B='West Winter - West Spring'; D=-11.4150487495947; E=2.06160960069476; G=-5.53695944457568; H=0.00000217297647409076},
    @{B='West Winter - West Summer'; D=-24.6820256807924; E=1.94332611709243; G=-12.7009180104682; H=0}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r++
}
Write-Host "Applied $($rows.Count) rows"
